$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.348, -0.044, 0.06900000000000001, 75, -8, 41),
    @(0.06, 0.224, 0.505, 65, -54, 90),
    @(0.09, -0.15, 0.513, 105, 74, 43),
    @(0.016, 0.381, -0.131, 145, -55, 65),
    @(-0.116, -0.176, -0.21, -126, 49, 113),
    @(-0.046, -0.136, 0.443, -46, -11, 44),
    @(0.09, -0.493, 0.331, -46, -70, -113),
    @(-0.129, -0.012, 0.098, -178, -3, -153),
    @(-0.124, 0.095, 0.171, -173, -4, -9),
    @(-0.023, -0.027, 0.363, -89, 70, -36)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}
